$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  C = 1.01716982815488;   D = 0.3162564994931649 },
    @{ Row = 3;  C = 1.089166687490565;  D = 0.2837403477047964 },
    @{ Row = 4;  C = 0.4291860951955537; D = 0.6704952402572961 },
    @{ Row = 5;  C = -0.3277491580830804; D = 0.7451110494170994 },
    @{ Row = 6;  C = 0.04203661712253944; D = 0.9667154040331787 },
    @{ Row = 7;  C = -0.8199806614268734; D = 0.4179386983821669 },
    @{ Row = 8;  C = -1.046847789468254;  D = 0.302555165211988 },
    @{ Row = 9;  C = -0.9444565343425877; D = 0.351600056886006 },
    @{ Row = 10; C = -1.176259306712132;  D = 0.2476615635043664 },
    @{ Row = 11; C = -0.6106451015356913; D = 0.5454952792654466 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    $ws.Cells.Item($u.Row, 4).Value = $u.D
}
